$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 12
$ws.Range("H12").Value = 965.6667
$ws.Range("I12").Value = 974
$ws.Range("J12").Value = 949
$ws.Range("K12").Value = 974
$ws.Range("L12").Value = 949
$ws.Range("M12").Value = -804
$ws.Range("N12").Value = -1289

# Row 58
$ws.Range("H58").Value = 4448.846
$ws.Range("I58").Value = 67.5
$ws.Range("J58").Value = 5245.4546
$ws.Range("K58").Value = 202.5
$ws.Range("L58").Value = 15736.3638
$ws.Range("M58").Value = -52.5
$ws.Range("N58").Value = -16036.3638

# Row 76
$ws.Range("H76").Value = 6375.9443
$ws.Range("I76").Value = 5639.1665
$ws.Range("K76").Value = 5639.1665
$ws.Range("M76").Value = -5324.1665

# Row 79
$ws.Range("H79").Value = 6375.9443
$ws.Range("I79").Value = 5639.1665
$ws.Range("K79").Value = 5639.1665
$ws.Range("M79").Value = -4547.1665

# Row 131
$ws.Range("H131").Value = 6153.5454

# Row 132
$ws.Range("H132").Value = 1392.0769
$ws.Range("I132").Value = 1347.619
$ws.Range("K132").Value = 4042.857
$ws.Range("M132").Value = -1512.857

# Row 138
$ws.Range("H138").Value = 5163.018
$ws.Range("I138").Value = 2034
$ws.Range("J138").Value = 5475.92
$ws.Range("K138").Value = 6102
$ws.Range("L138").Value = 16427.76
$ws.Range("M138").Value = -962
$ws.Range("N138").Value = -26707.76

$ws = $wb.Worksheets.Item("ARM")
# Row 8
$ws.Range("H8").Value = 600961.2
$ws.Range("I8").Value = 1500500
$ws.Range("J8").Value = 1268.6666
$ws.Range("K8").Value = 1500500
$ws.Range("L8").Value = 1268.6666
$ws.Range("M8").Value = -1500356
$ws.Range("N8").Value = -1556.6666

# Row 10
$ws.Range("H10").Value = 805
$ws.Range("I10").Value = 0
$ws.Range("K10").Value = 0
$ws.Range("M10").ClearContents()

# Row 12
$ws.Range("H12").Value = 300.66666
$ws.Range("I12").Value = 401
$ws.Range("K12").Value = 401
$ws.Range("M12").Value = -228

# Row 14
$ws.Range("H14").Value = 2069
$ws.Range("I14").Value = 2650
$ws.Range("J14").Value = 907
$ws.Range("K14").Value = 2650
$ws.Range("L14").Value = 907
$ws.Range("M14").Value = -2475
$ws.Range("N14").Value = -1257

# Row 74
$ws.Range("H74").Value = 4700.7144
$ws.Range("I74").Value = 1596
$ws.Range("K74").Value = 1596
$ws.Range("M74").Value = -722

# Row 77
$ws.Range("H77").Value = 4700.7144
$ws.Range("I77").Value = 1596
$ws.Range("K77").Value = 7980
$ws.Range("M77").Value = -3612

# Row 110
$ws.Range("H110").Value = 5383.857
$ws.Range("I110").Value = 5383.857
$ws.Range("K110").Value = 5383.857
$ws.Range("M110").Value = -3338.857

# Row 132
$ws.Range("H132").Value = 1729.3235
$ws.Range("I132").Value = 1638.871
$ws.Range("K132").Value = 4916.613
$ws.Range("M132").Value = -2386.613

$ws = $wb.Worksheets.Item("BSM")
# Row 86
$ws.Range("H86").Value = 2751.375
$ws.Range("I86").Value = 1241
$ws.Range("J86").Value = 5268.6665
$ws.Range("K86").Value = 1241
$ws.Range("L86").Value = 5268.6665
$ws.Range("M86").Value = -118
$ws.Range("N86").Value = -7514.6665

# Row 89
$ws.Range("H89").Value = 2751.375
$ws.Range("I89").Value = 1241
$ws.Range("J89").Value = 5268.6665
$ws.Range("K89").Value = 6205
$ws.Range("L89").Value = 26343.3325
$ws.Range("M89").Value = -589
$ws.Range("N89").Value = -37575.3325

# Row 105
$ws.Range("H105").Value = 3932.9062
$ws.Range("I105").Value = 3159.7144
$ws.Range("K105").Value = 3159.7144
$ws.Range("M105").Value = -1412.7144

$ws = $wb.Worksheets.Item("CRP")
# Row 58
$ws.Range("H58").Value = 4033.9333
$ws.Range("I58").Value = 2854.6667
$ws.Range("K58").Value = 2854.6667
$ws.Range("M58").Value = -2651.6667

# Row 99
$ws.Range("H99").Value = 16838.059
$ws.Range("I99").Value = 13981.714
$ws.Range("K99").Value = 13981.714
$ws.Range("M99").Value = -12483.714

# Row 126
$ws.Range("H126").Value = 16838.059
$ws.Range("I126").Value = 13981.714
$ws.Range("K126").Value = 41945.142
$ws.Range("M126").Value = -39475.142

# Row 134
$ws.Range("H134").Value = 2745.8823
$ws.Range("I134").Value = 1376.4166
$ws.Range("K134").Value = 4129.2498
$ws.Range("M134").Value = -1594.2498

# Row 136
$ws.Range("H136").Value = 4033.9333
$ws.Range("I136").Value = 2854.6667
$ws.Range("K136").Value = 8564.000100000001
$ws.Range("M136").Value = -6014.000100000001

$ws = $wb.Worksheets.Item("CUL")
# Row 7
$ws.Range("H7").Value = 6666827.5
$ws.Range("I7").Value = 10000077
$ws.Range("K7").Value = 30000231
$ws.Range("M7").Value = -30000119

# Row 34
$ws.Range("H34").Value = 1554.75
$ws.Range("I34").Value = 1139
$ws.Range("J34").Value = 1831.9166
$ws.Range("K34").Value = 3417
$ws.Range("L34").Value = 5495.7498
$ws.Range("M34").Value = -3333
$ws.Range("N34").Value = -5663.7498

# Row 39
$ws.Range("H39").Value = 2554.889
$ws.Range("J39").Value = 3502
$ws.Range("L39").Value = 10506
$ws.Range("N39").Value = -11094

# Row 55
$ws.Range("H55").Value = 1000000
$ws.Range("I55").Value = 1000000
$ws.Range("J55").Value = 0
$ws.Range("K55").Value = 3000000
$ws.Range("L55").Value = 0
$ws.Range("M55").Value = -2999823
$ws.Range("N55").ClearContents()

# Row 75
$ws.Range("H75").Value = 398.66666
$ws.Range("J75").Value = 250
$ws.Range("L75").Value = 750
$ws.Range("N75").Value = -2746

# Row 78
$ws.Range("H78").Value = 398.66666
$ws.Range("J78").Value = 250
$ws.Range("L78").Value = 2250
$ws.Range("N78").Value = -12234

# Row 86
$ws.Range("H86").Value = 294.83334
$ws.Range("J86").Value = 272.25
$ws.Range("L86").Value = 816.75
$ws.Range("N86").Value = -3188.75

# Row 89
$ws.Range("H89").Value = 294.83334
$ws.Range("J89").Value = 272.25
$ws.Range("L89").Value = 2450.25
$ws.Range("N89").Value = -14306.25

# Row 94
$ws.Range("H94").Value = 0
$ws.Range("I94").Value = 0
$ws.Range("K94").Value = 0
$ws.Range("M94").ClearContents()

# Row 103
$ws.Range("H103").Value = 384.2857
$ws.Range("I103").Value = 312.5
$ws.Range("J103").Value = 480
$ws.Range("K103").Value = 937.5
$ws.Range("L103").Value = 1440
$ws.Range("M103").Value = -58.5
$ws.Range("N103").Value = -3198

# Row 121
$ws.Range("H121").Value = 1170
$ws.Range("I121").Value = 212.5
$ws.Range("K121").Value = 637.5
$ws.Range("M121").Value = 672.5

# Row 131
$ws.Range("H131").Value = 3276.6453
$ws.Range("J131").Value = 3651.5
$ws.Range("L131").Value = 10954.5
$ws.Range("N131").Value = -21034.5

$ws = $wb.Worksheets.Item("GSM")
# Row 33
$ws.Range("H33").Value = 19998
$ws.Range("J33").Value = 19998
$ws.Range("L33").Value = 19998
$ws.Range("N33").Value = -20502

# Row 36
$ws.Range("H36").Value = 3000
$ws.Range("J36").Value = 3000
$ws.Range("L36").Value = 3000
$ws.Range("N36").Value = -3970

# Row 132
$ws.Range("H132").Value = 2469.077
$ws.Range("I132").Value = 1938.909
$ws.Range("J132").Value = 5385
$ws.Range("K132").Value = 5816.727000000001
$ws.Range("L132").Value = 16155
$ws.Range("M132").Value = -3286.727000000001
$ws.Range("N132").Value = -21215

# Row 136
$ws.Range("H136").Value = 25966.215
$ws.Range("J136").Value = 25966.215
$ws.Range("L136").Value = 77898.645
$ws.Range("N136").Value = -82998.645

$ws = $wb.Worksheets.Item("LTW")
# Row 100
$ws.Range("H100").Value = 4370.1
$ws.Range("I100").Value = 2077.889
$ws.Range("K100").Value = 2077.889
$ws.Range("M100").Value = -1536.889

# Row 132
$ws.Range("H132").Value = 5154.5557
$ws.Range("J132").Value = 5096
$ws.Range("L132").Value = 15288
$ws.Range("N132").Value = -20348

$ws = $wb.Worksheets.Item("WVR")
# Row 15
$ws.Range("H15").Value = 20999.666
$ws.Range("I15").Value = 20999
$ws.Range("K15").Value = 20999
$ws.Range("M15").Value = -20711

# Row 113
$ws.Range("H113").Value = 1465.1143
$ws.Range("I113").Value = 1317
$ws.Range("J113").Value = 1641
$ws.Range("K113").Value = 3951
$ws.Range("L113").Value = 4923
$ws.Range("M113").Value = -1781
$ws.Range("N113").Value = -9263

# Row 122
$ws.Range("H122").Value = 1826
$ws.Range("I122").Value = 1826
$ws.Range("K122").Value = 5478
$ws.Range("M122").Value = -3028
